$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting from the last two existing data rows (4 and 5) onto
# the two new rows (6 and 7) so the new rows inherit the same styles.
$ws.Range("A4:M4").Copy()
$ws.Range("A6:M6").PasteSpecial(-4122)
$ws.Range("A5:M5").Copy()
$ws.Range("A7:M7").PasteSpecial(-4122)

# New shared strings need to be introduced in the same order as the
# original edit so the rebuilt sharedStrings table matches exactly:
# TubeQPCRResults0005, TubeQPCRResults0006, UnknownPlate, A2
$ws.Range("A6").Value = "TubeQPCRResults0005"
$ws.Range("A7").Value = "TubeQPCRResults0006"
$ws.Range("E6").Value = "UnknownPlate"
$ws.Range("C7").Value = "A2"

# Row 6 - new "UnknownPlate" well A1 result
$ws.Range("B6").Value = "Non-Negative"
$ws.Range("C6").Value = "A1"
$ws.Range("D6").Value = "20200715_PVE690RLR_1"
$ws.Range("F6").Value = "Rack001"
$ws.Range("G6").Value = "B02"
$ws.Range("H6").Value = "Undetermined"
$ws.Range("I6").Value = 200
$ws.Range("J6").Value = 19.6210052131733
$ws.Range("K6").Value = 400
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 1.96854084629378

# Row 7 - new "UnknownPlate" well A2 result
$ws.Range("B7").Value = "Negative"
$ws.Range("D7").Value = "20200715_PVE690RLR_1"
$ws.Range("E7").Value = "UnknownPlate"
$ws.Range("F7").Value = "Rack001"
$ws.Range("G7").Value = "B02"
$ws.Range("H7").Value = "Undetermined"
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 21.98765432
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 0

$ws.Range("E7").Select()
